# "CSV angepasst, neues File"
# Append 7 new log rows (217-223) to the "created leaders" log sheet,
# mirroring the existing row layout:
#   A: timestamp, B: document, C: leader_guid, D: type, E: dimstyle, F: csv

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2025-10-16T10:33:57", "Untitled", "efa63824-f019-405b-b75a-639ac869047c", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:39:26", "Untitled", "94f76575-f7e5-441e-be91-ef4476f6181d", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:42:11", "Untitled", "c0d59d16-dda2-4acb-b083-f2b17560a1b2", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:42:50", "Untitled", "6894475a-8d1e-4644-a7eb-8938855707da", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:43:19", "Untitled", "f3c95967-ae33-496d-b12d-497bebd44286", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:43:34", "Untitled", "726ab9d8-5e81-44bd-bd0a-6e489a846532", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv"),
    @("2025-10-16T10:43:52", "Untitled", "c881a14f-457d-43ea-bf8e-505042d6be41", "rahmentuere", "Standard 1:10 Rahmenbeschriftung", "rahmentuereBlockrahmen.csv")
)

$startRow = 217
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $entry = $newRows[$i]
    $row = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $entry[$col - 1]
        # keep the same default ("Normal") vertical alignment the rest of
        # the log uses instead of inheriting the blank row's bottom-aligned
        # row-level style override
        $cell.VerticalAlignment = -4107
    }
}
